# TrackingOrder test-result refresh (headless-browser run):
# the MDSi test run appended a new batch of job numbers to the shared
# string table and re-pointed the three result rows (B2:B4) at the
# latest three job numbers for this pass.
#
# Values like "32381731" look numeric, so a plain Range.Value assignment
# would store them as numbers. Mark the cells as Text first (matches how
# the existing B column values are stored as shared strings), write the
# values, then drop the cells back to the default "Normal" style so the
# formatting matches the rest of column B (no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:B4")
$rng.NumberFormat = "@"

$ws.Range("B2").Value = "32381731"
$ws.Range("B3").Value = "32381732"
$ws.Range("B4").Value = "32381733"

$rng.Style = "Normal"
